function Get-ParaIndexByText($doc, $targetText) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        $ptext = $p.Range.Text
        $ptext = $ptext.TrimEnd([char]13, [char]7)
        if ($ptext -eq $targetText) {
            return $i
        }
    }
    return -1
}

function Insert-ItalicParagraphAfter($doc, $afterIndex, $text) {
    $pAfter = $doc.Paragraphs($afterIndex)
    $pAfter.Range.InsertParagraphAfter()
    $pNew = $doc.Paragraphs($afterIndex + 1)
    $pNew.Range.Text = $text
    $rOnly = $doc.Range($pNew.Range.Start, $pNew.Range.End - 1)
    $rOnly.Font.Italic = 1
}

$d = $word.ActiveDocument

# 1) Creditos-trabalho 2 -> 1
$old1 = @"
Créditos-trabalho: 2
"@
$new1 = @"
Créditos-trabalho: 1
"@
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) Carga horaria 75h -> 45h
$old2 = @"
Carga horária: 75 h
"@
$new2 = @"
Carga horária: 45 h
"@
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Ativacao date
$old3 = @"
Ativação: 01/01/2012
"@
$new3 = @"
Ativação: 01/01/2022
"@
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4) Objetivos paragraph text replace
$old4 = @"
O TG tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos previamente adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro ambiental. O TG deverá evidenciar a atualização de conhecimentos, a objetividade e a reflexão pessoal. No TG II, o aluno deverá escrever uma monografia sobre o plano de trabalho aprovado em TG I e apresentá-la a uma banca examinadora.
"@
$new4 = @"
O Trabalho de Graduação 2 tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O desenvolvimento do Trabalho de Graduação 2 deverá seguir o modelo escolhido no Trabalho de Graduação 1, e permitirá o uso de tecnologias digitais ou outras metodologias para desenvolvimento de conteúdo ou produto aplicável que utilize conteúdo da Engenharia Ambiental ou correlatas.
"@
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 7) Programa resumido paragraph text replace
$old7 = @"
Elaborar a monografia de Trabalho de Graduação e apresentá-la perante uma banca de examinadores.
"@
$new7 = @"
O aluno deverá desenvolver o Trabalho de Graduação 2 seguindo o modelo escolhido e já utilizado no Trabalho de Graduação 1. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados a Engenharia Ambiental, e ao final, o Trabalho de Graduação deverá ser avaliado por banca avaliadora.
"@
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# 9) Programa paragraph text replace
$old9 = @"
O programa da disciplina será constituído pelas seguintes etapas: 1. Elaborar a monografia de acordo com a proposta aprovada na disciplina TG I. 2. Definição e divulgação da data de apresentação após a entrega da monografia com antecedência de, no mínimo, 15 dias úteis. 3. Definição da banca de examinadores, sendo constituída pelo professor orientador e por dois professores convidados. 4. Apresentação e avaliação do TG. 5. Divulgação da avaliação. 6. Em caso de aprovação, entrega do exemplar final da monografia (cópia impressa e eletrônica). 7. Envio da nota final.
"@
$new9 = @"
Para ambos os modelos (artigo ou produto): O aluno deverá dar continuidade ao desenvolvimento do Trabalho iniciado no Trabalho de Graduação 1. O programa da disciplina é constituído pelas seguintes etapas: 1) Desenvolvimento do tema com base nas atividades de cronograma aprovados no Trabalho de Graduação 1. 2) Desenvolvimento do texto final, conforme modelo fornecido pelos responsáveis da disciplina. 3) Entrega da versão final do texto, com aprovação do orientador e com a indicação da banca de avaliação. 4) Avaliação e atribuição de notas pela banca avaliadora, que pode ou não ser a mesma composta anteriormente no Trabalho de Graduação 1, a critério do orientador.
"@
$d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2) | Out-Null

# 11) Metodo text replace
$old11 = @"
Em função da natureza deste curso, a avaliação será feita pela elaboração e apresentação de uma monografia.
"@
$new11 = @"
Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)
"@
$d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2) | Out-Null

# 12) Criterio text replace
$old12 = @"
Avaliação e atribuição de nota do Trabalho de Graduação por uma banca examinadora.
"@
$new12 = @"
Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final.Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil.
"@
$d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2) | Out-Null

# 13) Bibliografia text replace
$old13 = @"
Definida na monografia do TG.
"@
$new13 = @"
A ser definido no decorrer de cada projeto
"@
$d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, 1, $false, $new13, 2) | Out-Null

# 6) Add second docente run with line break within Docentes paragraph
$docenteIdx = Get-ParaIndexByText $d "4780627 - Ana Lucia Gabas Ferreira"
$pDoc = $d.Paragraphs($docenteIdx)
$lineBreak = [char]11
$pDoc.Range.InsertAfter($lineBreak + "7455355 - Robson da Silva Rocha")

# 5) Insert new English italic paragraph after Objetivos paragraph (new PT text)
$en5 = @"
Graduation Work 2 aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to carry out activities, which are part of the professional performance profile of the environmental engineer. The development of Graduate Work 2 should follow the model chosen in Graduate Work 1, and will allow the use of digital technologies or other methodologies for the development of content or applicable product that uses Environmental Engineering or related content.
"@
$objIdx = Get-ParaIndexByText $d $new4
Insert-ItalicParagraphAfter $d $objIdx $en5

# 8) Insert new English italic paragraph after Programa resumido paragraph (new PT text)
$en8 = @"
The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel.
"@
$resIdx = Get-ParaIndexByText $d $new7
Insert-ItalicParagraphAfter $d $resIdx $en8

# 10) Insert new English italic paragraph after Programa paragraph (new PT text)
$en10 = @"
For both models (article or product): The student must continue with the development of the Work initiated in the Graduate Work 1. The course program consists of the following steps: 1) Development of the theme based on the schedule activities approved in the Graduate Work 1. 2) Development of the final text, according to the model provided by those responsible for the discipline. 3) Delivery of the final version of the text, with the approval of the advisor and with the indication of the evaluation board. 4) Evaluation and grading by the examining board, which may or may not be the same previously composed in Graduate Work 1, at the discretion of the advisor.
"@
$progIdx = Get-ParaIndexByText $d $new9
Insert-ItalicParagraphAfter $d $progIdx $en10

Write-Output "edit complete"
